# Auto-generated: updates Leve profit calculation columns (H-N) across all
# sheets to reflect refreshed market-board pricing, per scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 426.6842
$ws.Range("I33").Value = 403.0625
$ws.Range("K33").Value = 403.0625
$ws.Range("M33").Value = -174.0625
$ws.Range("H51").Value = 2072.6
$ws.Range("I51").Value = 1430
$ws.Range("J51").Value = 2233.25
$ws.Range("K51").Value = 1430
$ws.Range("L51").Value = 2233.25
$ws.Range("M51").Value = -946
$ws.Range("N51").Value = -3201.25
$ws.Range("H92").Value = 1608.4231
$ws.Range("I92").Value = 1665.8235
$ws.Range("K92").Value = 1665.8235
$ws.Range("M92").Value = -417.8235
$ws.Range("H132").Value = 7581819.5
$ws.Range("I132").Value = 10106264
$ws.Range("J132").Value = 8487.637000000001
$ws.Range("K132").Value = 30318792
$ws.Range("L132").Value = 25462.911
$ws.Range("M132").Value = -30316262
$ws.Range("N132").Value = -30522.911

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4319.0137
$ws.Range("I32").Value = 3826.877
$ws.Range("K32").Value = 3826.877
$ws.Range("M32").Value = -3539.877
$ws.Range("H61").Value = 90910990
$ws.Range("I61").Value = 111112770
$ws.Range("K61").Value = 111112770
$ws.Range("M61").Value = -111112558
$ws.Range("H132").Value = 1367.9736
$ws.Range("I132").Value = 1278.5555
$ws.Range("K132").Value = 3835.6665
$ws.Range("M132").Value = -1305.6665
$ws.Range("H136").Value = 90910990
$ws.Range("I136").Value = 111112770
$ws.Range("K136").Value = 333338310
$ws.Range("M136").Value = -333335760

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1396.9744
$ws.Range("I20").Value = 1052.1923
$ws.Range("J20").Value = 2086.5386
$ws.Range("K20").Value = 1052.1923
$ws.Range("L20").Value = 2086.5386
$ws.Range("M20").Value = -805.1922999999999
$ws.Range("N20").Value = -2580.5386
$ws.Range("H86").Value = 3041.2354
$ws.Range("I86").Value = 3112.5625
$ws.Range("K86").Value = 3112.5625
$ws.Range("M86").Value = -1989.5625
$ws.Range("H89").Value = 3041.2354
$ws.Range("I89").Value = 3112.5625
$ws.Range("K89").Value = 15562.8125
$ws.Range("M89").Value = -9946.8125
$ws.Range("H132").Value = 57999.75
$ws.Range("J132").Value = 57999.75
$ws.Range("L132").Value = 57999.75
$ws.Range("N132").Value = -68119.75
$ws.Range("H134").Value = 3392.8667
$ws.Range("I134").Value = 1005.4167
$ws.Range("J134").Value = 12942.667
$ws.Range("K134").Value = 3016.2501
$ws.Range("L134").Value = 38828.001
$ws.Range("M134").Value = -481.2501000000002
$ws.Range("N134").Value = -43898.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1927.9688
$ws.Range("I31").Value = 1840.5927
$ws.Range("J31").Value = 2399.8
$ws.Range("K31").Value = 1840.5927
$ws.Range("L31").Value = 2399.8
$ws.Range("M31").Value = -1545.5927
$ws.Range("N31").Value = -2989.8
$ws.Range("H34").Value = 1927.9688
$ws.Range("I34").Value = 1840.5927
$ws.Range("J34").Value = 2399.8
$ws.Range("K34").Value = 1840.5927
$ws.Range("L34").Value = 2399.8
$ws.Range("M34").Value = -1638.5927
$ws.Range("N34").Value = -2803.8
$ws.Range("H58").Value = 893.12195
$ws.Range("I58").Value = 828.8788
$ws.Range("J58").Value = 1158.125
$ws.Range("K58").Value = 828.8788
$ws.Range("L58").Value = 1158.125
$ws.Range("M58").Value = -625.8788
$ws.Range("N58").Value = -1564.125
$ws.Range("H86").Value = 1812609
$ws.Range("I86").Value = 2784204
$ws.Range("K86").Value = 2784204
$ws.Range("M86").Value = -2783081
$ws.Range("H89").Value = 1812609
$ws.Range("I89").Value = 2784204
$ws.Range("K89").Value = 13921020
$ws.Range("M89").Value = -13915404
$ws.Range("H105").Value = 905
$ws.Range("I105").Value = 905
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 905
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 842
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 578.7727
$ws.Range("J107").Value = 708.6667
$ws.Range("L107").Value = 708.6667
$ws.Range("N107").Value = -4548.6667
$ws.Range("H132").Value = 3150.1206
$ws.Range("I132").Value = 3429.58
$ws.Range("J132").Value = 1403.5
$ws.Range("K132").Value = 10288.74
$ws.Range("L132").Value = 4210.5
$ws.Range("M132").Value = -7758.74
$ws.Range("N132").Value = -9270.5
$ws.Range("H134").Value = 11629232
$ws.Range("I134").Value = 1370.6316
$ws.Range("K134").Value = 4111.8948
$ws.Range("M134").Value = -1576.8948
$ws.Range("H136").Value = 893.12195
$ws.Range("I136").Value = 828.8788
$ws.Range("J136").Value = 1158.125
$ws.Range("K136").Value = 2486.6364
$ws.Range("L136").Value = 3474.375
$ws.Range("M136").Value = 63.36360000000013
$ws.Range("N136").Value = -8574.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2709.7083
$ws.Range("J39").Value = 2631.65
$ws.Range("L39").Value = 7894.950000000001
$ws.Range("N39").Value = -8482.950000000001
$ws.Range("H107").Value = 4887.6816
$ws.Range("I107").Value = 341.7143
$ws.Range("J107").Value = 7009.1333
$ws.Range("K107").Value = 1025.1429
$ws.Range("L107").Value = 21027.3999
$ws.Range("M107").Value = 894.8571000000002
$ws.Range("N107").Value = -24867.3999
$ws.Range("H113").Value = 729.56525
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 729.56525
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2188.69575
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6528.69575
$ws.Range("H131").Value = 45456350
$ws.Range("J131").Value = 2512.3572
$ws.Range("L131").Value = 7537.071599999999
$ws.Range("N131").Value = -17617.0716
$ws.Range("H132").Value = 1877.1538
$ws.Range("J132").Value = 1937.875
$ws.Range("L132").Value = 17440.875
$ws.Range("N132").Value = -22500.875
$ws.Range("H139").Value = 2007.55
$ws.Range("I139").Value = 2340.4
$ws.Range("J139").Value = 1674.7
$ws.Range("K139").Value = 7021.200000000001
$ws.Range("L139").Value = 5024.1
$ws.Range("M139").Value = -1881.200000000001
$ws.Range("N139").Value = -15304.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 35717890
$ws.Range("I102").Value = 50001600
$ws.Range("J102").Value = 8607
$ws.Range("K102").Value = 50001600
$ws.Range("L102").Value = 8607
$ws.Range("M102").Value = -49999978
$ws.Range("N102").Value = -11851
$ws.Range("H126").Value = 2832.6667
$ws.Range("I126").Value = 1850
$ws.Range("J126").Value = 3487.7778
$ws.Range("K126").Value = 5550
$ws.Range("L126").Value = 10463.3334
$ws.Range("M126").Value = -3080
$ws.Range("N126").Value = -15403.3334
$ws.Range("H132").Value = 1958.3928
$ws.Range("I132").Value = 1601.5217
$ws.Range("K132").Value = 4804.5651
$ws.Range("M132").Value = -2274.5651

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 900
$ws.Range("I30").Value = 900
$ws.Range("K30").Value = 900
$ws.Range("M30").Value = -792
$ws.Range("H40").Value = 2451.6155
$ws.Range("I40").Value = 2324.6365
$ws.Range("J40").Value = 3150
$ws.Range("K40").Value = 2324.6365
$ws.Range("M40").Value = -2188.6365
$ws.Range("N40").Value = -3422
$ws.Range("H55").Value = 392.79166
$ws.Range("I55").Value = 242.75
$ws.Range("K55").Value = 242.75
$ws.Range("M55").Value = -69.75
$ws.Range("H132").Value = 27405.896
$ws.Range("I132").Value = 1539.5
$ws.Range("K132").Value = 4618.5
$ws.Range("M132").Value = -2088.5
$ws.Range("H136").Value = 4986.448
$ws.Range("I136").Value = 5304.28
$ws.Range("K136").Value = 15912.84
$ws.Range("M136").Value = -13362.84

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 9999
$ws.Range("J25").Value = 9999
$ws.Range("L25").Value = 9999
$ws.Range("N25").Value = -10585
$ws.Range("H113").Value = 291.68
$ws.Range("I113").Value = 181.3
$ws.Range("J113").Value = 365.26666
$ws.Range("K113").Value = 543.9000000000001
$ws.Range("L113").Value = 1095.79998
$ws.Range("M113").Value = 1626.1
$ws.Range("N113").Value = -5435.79998
$ws.Range("H132").Value = 2452.0344
$ws.Range("I132").Value = 3292.8823
$ws.Range("J132").Value = 1260.8334
$ws.Range("K132").Value = 9878.6469
$ws.Range("L132").Value = 3782.5002
$ws.Range("M132").Value = -7348.6469
$ws.Range("N132").Value = -8842.5002
$ws.Range("H136").Value = 564.0625
$ws.Range("I136").Value = 364.81482
$ws.Range("J136").Value = 1640
$ws.Range("K136").Value = 1094.44446
$ws.Range("L136").Value = 4920
$ws.Range("M136").Value = 1455.55554
$ws.Range("N136").Value = -10020
